$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$ws.Range("A1").Value = 45311
$ws.Range("D33").Value = 26.828
$ws.Range("D34").Value = 19.256
